# #914: Cleanup node SshProxy initialization and relocate home folders
#
# Three textual tweaks inside "Ubuntu-20.04 XenServer Template.docx":
#   1. & 2. The two "neon prepare node-template <EN-DASH>xenserver" code
#      samples had an autocorrected en-dash (\u2013) instead of a literal
#      double-hyphen ("--") in front of "xenserver". Restore "--".
#   3. The "Delete the template from XenServer." bullet gains a
#      parenthetical "(if you wish)" before the final period, splitting
#      the former single run into three runs:
#      "from XenServer" / " (if you wish)" / ".".

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Fix 1 & 2: replace the leading en-dash (U+2013) in front of
# "xenserver" with a plain double-hyphen, without disturbing the
# surrounding runs/text.
# ---------------------------------------------------------------------
function Repair-DashXenServer {
    $rng = $d.Content
    $found = $rng.Find.Execute([string][char]0x2013 + "xenserver", $true, $false, $false, `
                                $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        return $false
    }

    # Narrow range covering just the single en-dash character.
    $dashRng = $d.Range($rng.Start, $rng.Start + 1)
    $dashRng.Text = "--"

    # Touch (and immediately revert) a character property on the now
    # two-character run so the engine keeps it as its own run instead of
    # silently re-merging it into the identically-formatted run that
    # precedes it ("neon prepare node-template ").
    $dashRng2 = $d.Range($rng.Start, $rng.Start + 2)
    $dashRng2.Bold = 1
    $dashRng2.Bold = 0

    return $true
}

# There are two occurrences in the document; fix both.
Repair-DashXenServer | Out-Null
Repair-DashXenServer | Out-Null

# ---------------------------------------------------------------------
# Fix 3: "from XenServer." -> "from XenServer" / " (if you wish)" / "."
# ---------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("from XenServer.", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)
if ($found) {
    # Collapse to the trailing period and insert the new clause right
    # before it, which naturally splits "from XenServer." into
    # "from XenServer" (untouched original run) + ". " (untouched
    # original run) with the newly inserted text in between.
    $rng.Collapse(0)
    $rng.MoveStart(1, -1) | Out-Null
    $rng.InsertBefore(" (if you wish)")

    $insertedLen = " (if you wish)".Length
    $insertedRng = $d.Range($rng.Start, $rng.Start + $insertedLen)

    # Same trick as above: force the inserted text to remain its own run
    # rather than being coalesced back into "from XenServer".
    $insertedRng.Bold = 1
    $insertedRng.Bold = 0
}
